# This workbook is a scraped product-price export. The crawl was re-run
# later the same day (new timestamp), and a handful of products that had
# gone out of online stock got an " - Online kein Bestand" suffix added
# to their ARIA-label text (column M), right before the trailing
# "<price> Schweizer Franken" portion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-12-30 12:55:45"
$newTimestamp = "2022-12-30 20:49:30"
$suffix = " - Online kein Bestand"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column letters: B=title/name, G=price, M=productAriaLabel, O=timestamp
$colB = 2
$colG = 7
$colM = 13
$colO = 15

# Row numbers (1-based, row 1 is the header) whose productAriaLabel (M)
# needs the "- Online kein Bestand" marker inserted before the price.
$outOfStockRows = @(3, 12, 15, 49)

for ($r = 2; $r -le $lastRow; $r++) {
    $tsCell = $ws.Cells.Item($r, $colO)
    if ($tsCell.Value2 -eq $oldTimestamp) {
        $tsCell.Value = $newTimestamp
    }
}

foreach ($r in $outOfStockRows) {
    $name = $ws.Cells.Item($r, $colB).Value2
    $price = $ws.Cells.Item($r, $colG).Value2
    $ws.Cells.Item($r, $colM).Value = $name + $suffix + " " + $price + " Schweizer Franken"
}
